$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1540
$ws.Range("I62").Value = 1540
$ws.Range("K62").Value = 1540
$ws.Range("M62").Value = -916

$ws.Range("H65").Value = 1540
$ws.Range("I65").Value = 1540
$ws.Range("K65").Value = 7700
$ws.Range("M65").Value = -4580

$ws.Range("H70").Value = 779
$ws.Range("I70").Value = 549
$ws.Range("J70").Value = 825
$ws.Range("K70").Value = 1647
$ws.Range("L70").Value = 2475
$ws.Range("M70").Value = -1377
$ws.Range("N70").Value = -3015

$ws.Range("H73").Value = 779
$ws.Range("I73").Value = 549
$ws.Range("J73").Value = 825
$ws.Range("K73").Value = 1647
$ws.Range("L73").Value = 2475
$ws.Range("M73").Value = -711
$ws.Range("N73").Value = -4347

$ws.Range("H106").Value = 3450
$ws.Range("I106").Value = 3299.6667
$ws.Range("J106").Value = 3514.4285
$ws.Range("K106").Value = 3299.6667
$ws.Range("L106").Value = 3514.4285
$ws.Range("M106").Value = -2668.6667
$ws.Range("N106").Value = -4776.4285

$ws.Range("H132").Value = 9094226
$ws.Range("I132").Value = 18183858
$ws.Range("J132").Value = 4594
$ws.Range("K132").Value = 54551574
$ws.Range("L132").Value = 13782
$ws.Range("M132").Value = -54549044
$ws.Range("N132").Value = -18842

$ws.Range("H137").Value = 1399.2927
$ws.Range("I137").Value = 1116.7059
$ws.Range("J137").Value = 2771.8572
$ws.Range("K137").Value = 3350.1177
$ws.Range("L137").Value = 8315.571599999999
$ws.Range("M137").Value = -800.1176999999998
$ws.Range("N137").Value = -13415.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5765.609
$ws.Range("I32").Value = 3467.3667
$ws.Range("K32").Value = 3467.3667
$ws.Range("M32").Value = -3180.3667

$ws.Range("H61").Value = 1667.8334
$ws.Range("I61").Value = 1500
$ws.Range("J61").Value = 1751.75
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 1751.75
$ws.Range("M61").Value = -1288
$ws.Range("N61").Value = -2175.75

$ws.Range("H74").Value = 38462970
$ws.Range("I74").Value = 55556620
$ws.Range("J74").Value = 2253.5
$ws.Range("K74").Value = 55556620
$ws.Range("L74").Value = 2253.5
$ws.Range("M74").Value = -55555746
$ws.Range("N74").Value = -4001.5

$ws.Range("H77").Value = 38462970
$ws.Range("I77").Value = 55556620
$ws.Range("J77").Value = 2253.5
$ws.Range("K77").Value = 277783100
$ws.Range("L77").Value = 11267.5
$ws.Range("M77").Value = -277778732
$ws.Range("N77").Value = -20003.5

$ws.Range("H97").Value = 678.8889
$ws.Range("I97").Value = 663.75
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 663.75
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -167.75
$ws.Range("N97").Value = -1792

$ws.Range("H110").Value = 2097
$ws.Range("I110").Value = 823
$ws.Range("J110").Value = 3583.3333
$ws.Range("K110").Value = 823
$ws.Range("L110").Value = 3583.3333
$ws.Range("M110").Value = 1222
$ws.Range("N110").Value = -7673.3333

$ws.Range("H122").Value = 1830.2565
$ws.Range("I122").Value = 1600.1154
$ws.Range("J122").Value = 2290.5386
$ws.Range("K122").Value = 4800.3462
$ws.Range("L122").Value = 6871.6158
$ws.Range("M122").Value = -2350.3462
$ws.Range("N122").Value = -11771.6158

$ws.Range("I127").Value = 30000
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 30000
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -25040
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 1645.5625
$ws.Range("I132").Value = 1268.5625
$ws.Range("J132").Value = 2022.5625
$ws.Range("K132").Value = 3805.6875
$ws.Range("L132").Value = 6067.6875
$ws.Range("M132").Value = -1275.6875
$ws.Range("N132").Value = -11127.6875

$ws.Range("H136").Value = 1667.8334
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 1751.75
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 5255.25
$ws.Range("M136").Value = -1950
$ws.Range("N136").Value = -10355.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13271.429
$ws.Range("J82").Value = 14600
$ws.Range("L82").Value = 14600
$ws.Range("N82").Value = -15366

$ws.Range("H85").Value = 13271.429
$ws.Range("J85").Value = 14600
$ws.Range("L85").Value = 14600
$ws.Range("N85").Value = -17252

$ws.Range("H94").Value = 1207.6666
$ws.Range("I94").Value = 526.875
$ws.Range("J94").Value = 1985.7142
$ws.Range("K94").Value = 526.875
$ws.Range("L94").Value = 1985.7142
$ws.Range("M94").Value = -75.875
$ws.Range("N94").Value = -2887.7142

$ws.Range("H105").Value = 1750220.9
$ws.Range("I105").Value = 4547394.5
$ws.Range("J105").Value = 1987.5
$ws.Range("K105").Value = 4547394.5
$ws.Range("L105").Value = 1987.5
$ws.Range("M105").Value = -4545647.5
$ws.Range("N105").Value = -5481.5

$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

$ws.Range("H134").Value = 11906008
$ws.Range("I134").Value = 13159120
$ws.Range("J134").Value = 1450
$ws.Range("K134").Value = 39477360
$ws.Range("L134").Value = 4350
$ws.Range("M134").Value = -39474825
$ws.Range("N134").Value = -9420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2750000
$ws.Range("I86").Value = 2750000
$ws.Range("K86").Value = 2750000
$ws.Range("M86").Value = -2748877

$ws.Range("H89").Value = 2750000
$ws.Range("I89").Value = 2750000
$ws.Range("K89").Value = 13750000
$ws.Range("M89").Value = -13744384

$ws.Range("H105").Value = 4134.7334
$ws.Range("I105").Value = 4601
$ws.Range("J105").Value = 3202.2
$ws.Range("K105").Value = 4601
$ws.Range("L105").Value = 3202.2
$ws.Range("M105").Value = -2854
$ws.Range("N105").Value = -6696.2

$ws.Range("H107").Value = 1225.9474
$ws.Range("I107").Value = 653.0714
$ws.Range("J107").Value = 2830
$ws.Range("K107").Value = 653.0714
$ws.Range("L107").Value = 2830
$ws.Range("M107").Value = 1266.9286
$ws.Range("N107").Value = -6670

$ws.Range("H132").Value = 2971.3076
$ws.Range("I132").Value = 2233.3333
$ws.Range("J132").Value = 4631.75
$ws.Range("K132").Value = 6699.999899999999
$ws.Range("L132").Value = 13895.25
$ws.Range("M132").Value = -4169.999899999999
$ws.Range("N132").Value = -18955.25

$ws.Range("H134").Value = 1109.6
$ws.Range("I134").Value = 1012
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 3036
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -501
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5452.2354
$ws.Range("I70").Value = 4726.1816
$ws.Range("J70").Value = 6783.3335
$ws.Range("K70").Value = 4726.1816
$ws.Range("L70").Value = 6783.3335
$ws.Range("M70").Value = -4456.1816
$ws.Range("N70").Value = -7323.3335

$ws.Range("H73").Value = 5452.2354
$ws.Range("I73").Value = 4726.1816
$ws.Range("J73").Value = 6783.3335
$ws.Range("K73").Value = 4726.1816
$ws.Range("L73").Value = 6783.3335
$ws.Range("M73").Value = -3790.1816
$ws.Range("N73").Value = -8655.333500000001

$ws.Range("H102").Value = 40429.69
$ws.Range("I102").Value = 51302.9
$ws.Range("J102").Value = 4185.6665
$ws.Range("K102").Value = 51302.9
$ws.Range("L102").Value = 4185.6665
$ws.Range("M102").Value = -49680.9
$ws.Range("N102").Value = -7429.6665

$ws.Range("H107").Value = 2052.5715
$ws.Range("I107").Value = 2751
$ws.Range("J107").Value = 1773.2
$ws.Range("K107").Value = 2751
$ws.Range("L107").Value = 1773.2
$ws.Range("M107").Value = -831
$ws.Range("N107").Value = -5613.2

$ws.Range("H132").Value = 8196.789000000001
$ws.Range("I132").Value = 12732.3
$ws.Range("J132").Value = 3157.3333
$ws.Range("K132").Value = 38196.89999999999
$ws.Range("L132").Value = 9471.999899999999
$ws.Range("M132").Value = -35666.89999999999
$ws.Range("N132").Value = -14531.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2347.96
$ws.Range("I61").Value = 1919.1666
$ws.Range("J61").Value = 2743.7693
$ws.Range("K61").Value = 1919.1666
$ws.Range("L61").Value = 2743.7693
$ws.Range("M61").Value = -1717.1666
$ws.Range("N61").Value = -3147.7693

$ws.Range("H100").Value = 2333.8333
$ws.Range("I100").Value = 2200.6
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2200.6
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1659.6
$ws.Range("N100").Value = -4082

$ws.Range("H113").Value = 2347.96
$ws.Range("I113").Value = 1919.1666
$ws.Range("J113").Value = 2743.7693
$ws.Range("K113").Value = 1919.1666
$ws.Range("L113").Value = 2743.7693
$ws.Range("M113").Value = 250.8334
$ws.Range("N113").Value = -7083.7693

$ws.Range("H122").Value = 6938.3447
$ws.Range("I122").Value = 7487.478
$ws.Range("J122").Value = 4833.3335
$ws.Range("K122").Value = 22462.434
$ws.Range("L122").Value = 14500.0005
$ws.Range("M122").Value = -20012.434
$ws.Range("N122").Value = -19400.0005

$ws.Range("H132").Value = 4032.1667
$ws.Range("I132").Value = 2627
$ws.Range("J132").Value = 4734.75
$ws.Range("K132").Value = 7881
$ws.Range("L132").Value = 14204.25
$ws.Range("M132").Value = -5351
$ws.Range("N132").Value = -19264.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1727.7778
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 1781.25
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 5343.75
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -10243.75

$ws.Range("H126").Value = 1170.5135
$ws.Range("I126").Value = 773.8182
$ws.Range("K126").Value = 2321.4546
$ws.Range("M126").Value = 148.5454

$ws.Range("H132").Value = 2022.381
$ws.Range("I132").Value = 1556.8572
$ws.Range("J132").Value = 2255.1428
$ws.Range("K132").Value = 4670.571599999999
$ws.Range("L132").Value = 6765.428400000001
$ws.Range("M132").Value = -2140.571599999999
$ws.Range("N132").Value = -11825.4284
